# Apply crypto price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "46.812.70"
$ws.Range("E2").Value = "  +6.41%  "

# Row 3
$ws.Range("D3").Value = "2.314.99"
$ws.Range("E3").Value = "  +5.39%  "

# Row 4
$ws.Range("E4").Value = "  -0.63%  "

# Row 5
$ws.Range("D5").Value = "'303.87"
$ws.Range("E5").Value = "  +2.92%  "

# Row 6
$ws.Range("D6").Value = "'103.41"
$ws.Range("E6").Value = "  +16.21%  "

# Row 7
$ws.Range("D7").Value = "'0.575"
$ws.Range("E7").Value = "  +1.86%  "

# Row 8
$ws.Range("E8").Value = "  -0.35%  "

# Row 9
$ws.Range("D9").Value = "'0.534"
$ws.Range("E9").Value = "  +10.57%  "

# Row 10
$ws.Range("D10").Value = "'37.39"
$ws.Range("E10").Value = "  +16.47%  "

# Row 11
$ws.Range("D11").Value = "'0.0808"
$ws.Range("E11").Value = "  +4.73%  "

# Row 12
$ws.Range("D12").Value = "'7.43"
$ws.Range("E12").Value = "  +9.61%  "

# Row 13
$ws.Range("E13").Value = "  +1.51%  "

# Row 14
$ws.Range("D14").Value = "2.669.23"
$ws.Range("E14").Value = "  +5.40%  "

# Row 15
$ws.Range("D15").Value = "2.313.43"
$ws.Range("E15").Value = "  +1.48%  "

# Row 16
$ws.Range("D16").Value = "'14.16"
$ws.Range("E16").Value = "  +8.17%  "

# Row 17
$ws.Range("D17").Value = "'0.830"
$ws.Range("E17").Value = "  +7.29%  "

# Row 18
$ws.Range("D18").Value = "46.834.44"
$ws.Range("E18").Value = "  +7.38%  "

# Row 19
$ws.Range("E19").Value = "  +23.78%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0953"
$ws.Range("E20").Value = "  +7.58%  "

# Row 21
$ws.Range("D21").Value = "'6.18"
$ws.Range("E21").Value = "  +5.88%  "

# Row 22
$ws.Range("D22").Value = "'67.23"
$ws.Range("E22").Value = "  +6.63%  "

# Row 23
$ws.Range("D23").Value = "'249.45"
$ws.Range("E23").Value = "  +6.41%  "

# Row 24
$ws.Range("D24").Value = "'2.98"
$ws.Range("E24").Value = "  +8.00%  "

# Row 25
$ws.Range("E25").Value = "  +8.81%  "

# Row 26
$ws.Range("E26").Value = "  -1.00%  "

# Row 27
$ws.Range("D27").Value = "'43.66"
$ws.Range("E27").Value = "  +20.83%  "

# Row 28
$ws.Range("E28").Value = "  +1.51%  "

# Row 29
$ws.Range("D29").Value = "'10.02"
$ws.Range("E29").Value = "  +8.35%  "

# Row 30
$ws.Range("D30").Value = "'20.26"
$ws.Range("E30").Value = "  +5.63%  "

# Row 31
$ws.Range("D31").Value = "'5.83"
$ws.Range("E31").Value = "  +10.48%  "

# Row 32
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'147.68"
$ws.Range("E32").Value = "  -0.35%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0807"
$ws.Range("E33").Value = "  +9.82%  "

# Row 34
$ws.Range("E34").Value = "  +4.80%  "

# Row 35
$ws.Range("D35").Value = "'3.14"
$ws.Range("E35").Value = "  +11.12%  "

# Row 36
$ws.Range("D36").Value = "'0.113"
$ws.Range("E36").Value = "  +9.97%  "

# Row 37
$ws.Range("E37").Value = "  +3.71%  "

# Row 38
$ws.Range("E38").Value = "  +10.24%  "

# Row 39
$ws.Range("D39").Value = "'15.98"
$ws.Range("E39").Value = "  +23.01%  "

# Row 40
$ws.Range("D40").Value = "'4.11"
$ws.Range("E40").Value = "  +16.65%  "

# Row 41
$ws.Range("D41").Value = "'3.50"
$ws.Range("E41").Value = "  +13.66%  "

# Row 42
$ws.Range("E42").Value = "  +9.42%  "

# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.66%  "

# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.99"
$ws.Range("E44").Value = "  +18.93%  "

# Row 45
$ws.Range("D45").Value = "1.853.52"
$ws.Range("E45").Value = "  +3.81%  "

# Row 46
$ws.Range("D46").Value = "'89.22"
$ws.Range("E46").Value = "  +22.09%  "

# Row 47
$ws.Range("D47").Value = "'0.199"
$ws.Range("E47").Value = "  +14.53%  "

# Row 48
$ws.Range("D48").Value = "'75.71"
$ws.Range("E48").Value = "  +16.88%  "

# Row 49
$ws.Range("E49").Value = "  +9.45%  "

# Row 50
$ws.Range("D50").Value = "'97.98"

# Row 51
$ws.Range("D51").Value = "'55.20"
$ws.Range("E51").Value = "  +12.47%  "
